$wb = $excel.ActiveWorkbook

# --- Summary sheet: update selection to C6 ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("C6").Select()

# --- Transactions sheet: update A2/A3 values and selection ---
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("A2").Value = 75
$wsTransactions.Range("A3").Value = 74
$wsTransactions.Activate()
$wsTransactions.Range("D3").Select()
